$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------
# Build the style "templates" we need on sheet2 by pulling formats
# from cells on Plan1 that already carry the right border, then
# resetting the font back to the default (size 11 / font 0) so the
# resulting cellXfs entries match fontId="0" + borderId="1".
# ---------------------------------------------------------------

# Style "center, no vertical" (same shape as existing xf idx 2) -> used
# for the two merged header/footer bars (A1:D1 and B5:E5).
$ws2.Range("A1:D1").Merge() | Out-Null
$ws1.Range("A1").Copy() | Out-Null
$ws2.Range("A1:D1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws2.Range("A1:D1").Font.Size = 11

$ws2.Range("B5:E5").Merge() | Out-Null
$ws1.Range("A1").Copy() | Out-Null
$ws2.Range("B5:E5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws2.Range("B5:E5").Font.Size = 11

# Style "center + vertical center" (same shape as existing xf idx 5) ->
# used for the two merged side columns (E1:E4 and A2:A5) plus the blank
# cells of column E / A that sit alongside the inner grid.
$ws2.Range("E1:E4").Merge() | Out-Null
$ws1.Range("A4").Copy() | Out-Null
$ws2.Range("E1:E4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws2.Range("E1:E4").Font.Size = 11

$ws2.Range("A2:A5").Merge() | Out-Null
$ws1.Range("A4").Copy() | Out-Null
$ws2.Range("A2:A5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws2.Range("A2:A5").Font.Size = 11

# Inner 3x3 grid (B2:D4) -> plain centered cells with the thin border on
# every side.
$ws1.Range("A1").Copy() | Out-Null
$ws2.Range("B2:D4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws2.Range("B2:D4").Font.Size = 11

# ---------------------------------------------------------------
# Fill in the cell values (shared-string order must match: A,B,C,D,E,
# F,G,H,I,J,L,M,N).
# ---------------------------------------------------------------
$ws2.Range("A1").Value = "A"
$ws2.Range("E1").Value = "B"

$ws2.Range("A2").Value = "C"
$ws2.Range("B2").Value = "D"
$ws2.Range("C2").Value = "E"
$ws2.Range("D2").Value = "F"

$ws2.Range("B3").Value = "G"
$ws2.Range("C3").Value = "H"
$ws2.Range("D3").Value = "I"

$ws2.Range("B4").Value = "J"
$ws2.Range("C4").Value = "L"
$ws2.Range("D4").Value = "M"

$ws2.Range("B5").Value = "N"

# ---------------------------------------------------------------
# Page setup + view state.
# ---------------------------------------------------------------
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

$ws2.Activate()
$ws2.Range("F3").Select()

Write-Host "done"
